$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet: insert a new blank column before column N (which
# held "Late"), pushing "Late" / "heading" / "Outstanding" one column to the
# right (N->O, O->P, P->Q). Give the freshly inserted column the same display
# width (11) as its left neighbour ("In Advance").
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = 10.166666666666666

# The active tab moves from "Edit Repayment Schedule" to "Repayment schedule",
# with cell K13 selected there; "Edit Repayment Schedule" keeps its previous
# selection (E14) but is no longer the active sheet.
$ws.Activate() | Out-Null
$ws.Range("K13").Select() | Out-Null
